$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New column header (F1) - shared string text, no special date style
$ws.Range("F1").Value = "11_03_2024"

# New column values (F2:F6)
$ws.Range("F2").Value = 3732
$ws.Range("F3").Value = 3189
$ws.Range("F4").Value = 4192
$ws.Range("F5").Value = 7163
$ws.Range("F6").Value = 363

# Update the active selection to match target (F7)
$ws.Range("F7").Select()
